$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = ""
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H113").Value = 3591.4546
$ws.Range("I113").Value = 3166.6667
$ws.Range("K113").Value = 3166.6667
$ws.Range("M113").Value = 87.33329999999978
$ws.Range("H132").Value = 14501394
$ws.Range("I132").Value = 15882093
$ws.Range("K132").Value = 47646279
$ws.Range("M132").Value = -47643749
$ws.Range("H133").Value = 38999
$ws.Range("J133").Value = 38999
$ws.Range("L133").Value = 38999
$ws.Range("N133").Value = -49119
$ws.Range("H138").Value = 1319.06
$ws.Range("I138").Value = 847.7560999999999
$ws.Range("J138").Value = 1646.5763
$ws.Range("K138").Value = 2543.2683
$ws.Range("L138").Value = 4939.7289
$ws.Range("M138").Value = 2596.7317
$ws.Range("N138").Value = -15219.7289
$ws.Range("H141").Value = 1355.5555
$ws.Range("I141").Value = 540
$ws.Range("K141").Value = 1620
$ws.Range("M141").Value = 3560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2231.889
$ws.Range("I61").Value = 1772
$ws.Range("J61").Value = 2599.8
$ws.Range("K61").Value = 1772
$ws.Range("L61").Value = 2599.8
$ws.Range("M61").Value = -1560
$ws.Range("N61").Value = -3023.8
$ws.Range("H102").Value = 33336138
$ws.Range("I102").Value = 33336138
$ws.Range("K102").Value = 33336138
$ws.Range("M102").Value = -33334516
$ws.Range("H110").Value = 1226.15
$ws.Range("I110").Value = 778.5
$ws.Range("K110").Value = 778.5
$ws.Range("M110").Value = 1266.5
$ws.Range("H136").Value = 2231.889
$ws.Range("I136").Value = 1772
$ws.Range("J136").Value = 2599.8
$ws.Range("K136").Value = 5316
$ws.Range("L136").Value = 7799.400000000001
$ws.Range("M136").Value = -2766
$ws.Range("N136").Value = -12899.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 252473920
$ws.Range("I105").Value = 252473920
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 252473920
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -252472173
$ws.Range("N105").Value = ""
$ws.Range("H134").Value = 9062.647000000001
$ws.Range("I134").Value = 1206.6
$ws.Range("K134").Value = 3619.8
$ws.Range("M134").Value = -1084.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1514.9584
$ws.Range("I31").Value = 1379.9546
$ws.Range("K31").Value = 1379.9546
$ws.Range("M31").Value = -1084.9546
$ws.Range("H34").Value = 1514.9584
$ws.Range("I34").Value = 1379.9546
$ws.Range("K34").Value = 1379.9546
$ws.Range("M34").Value = -1177.9546
$ws.Range("H105").Value = 730
$ws.Range("I105").Value = 662.5
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 662.5
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1084.5
$ws.Range("N105").Value = -4494
$ws.Range("H132").Value = 5746.6553
$ws.Range("I132").Value = 7787.0625
$ws.Range("J132").Value = 3235.3845
$ws.Range("K132").Value = 23361.1875
$ws.Range("L132").Value = 9706.1535
$ws.Range("M132").Value = -20831.1875
$ws.Range("N132").Value = -14766.1535
$ws.Range("H134").Value = 1739.9143
$ws.Range("I134").Value = 1725.3214
$ws.Range("K134").Value = 5175.9642
$ws.Range("M134").Value = -2640.9642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1238.4828
$ws.Range("I5").Value = 1247
$ws.Range("K5").Value = 3741
$ws.Range("M5").Value = -3629
$ws.Range("H86").Value = 850
$ws.Range("J86").Value = 850
$ws.Range("L86").Value = 2550
$ws.Range("N86").Value = -4922
$ws.Range("H89").Value = 850
$ws.Range("J89").Value = 850
$ws.Range("L89").Value = 7650
$ws.Range("N89").Value = -19506
$ws.Range("H122").Value = 698.75
$ws.Range("J122").Value = 799
$ws.Range("L122").Value = 7191
$ws.Range("N122").Value = -12091
$ws.Range("H131").Value = 16130255
$ws.Range("J131").Value = 1299.8246
$ws.Range("L131").Value = 3899.4738
$ws.Range("N131").Value = -13979.4738
$ws.Range("H132").Value = 923.75
$ws.Range("I132").Value = 440
$ws.Range("J132").Value = 1085
$ws.Range("K132").Value = 3960
$ws.Range("L132").Value = 9765
$ws.Range("M132").Value = -1430
$ws.Range("N132").Value = -14825
$ws.Range("H135").Value = 1238.4828
$ws.Range("I135").Value = 1247
$ws.Range("K135").Value = 11223
$ws.Range("M135").Value = -8688
$ws.Range("H139").Value = 3464.1875
$ws.Range("I139").Value = 3892.8
$ws.Range("J139").Value = 2749.8333
$ws.Range("K139").Value = 11678.4
$ws.Range("L139").Value = 8249.499899999999
$ws.Range("M139").Value = -6538.400000000001
$ws.Range("N139").Value = -18529.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1198.5
$ws.Range("I113").Value = 1198.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1198.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 971.5
$ws.Range("N113").Value = ""
$ws.Range("H126").Value = 2363.8462
$ws.Range("I126").Value = 2303.3333
$ws.Range("K126").Value = 6909.999899999999
$ws.Range("M126").Value = -4439.999899999999
$ws.Range("H132").Value = 2349.9524
$ws.Range("I132").Value = 1881.5333
$ws.Range("K132").Value = 5644.5999
$ws.Range("M132").Value = -3114.5999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3530
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3530
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3530
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -3802
$ws.Range("H55").Value = 221.16667
$ws.Range("I55").Value = 191.2
$ws.Range("J55").Value = 242.57143
$ws.Range("K55").Value = 191.2
$ws.Range("L55").Value = 242.57143
$ws.Range("M55").Value = -18.19999999999999
$ws.Range("N55").Value = -588.57143
$ws.Range("H112").Value = 46999.285
$ws.Range("J112").Value = 46999.285
$ws.Range("L112").Value = 46999.285
$ws.Range("N112").Value = -49953.285
$ws.Range("H122").Value = 35730144
$ws.Range("I122").Value = 35730144
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 107190432
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -107187982
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 21653.7
$ws.Range("I132").Value = 1295.9656
$ws.Range("K132").Value = 3887.8968
$ws.Range("M132").Value = -1357.8968

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 45457964
$ws.Range("I62").Value = 55558844
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 55558844
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -55558220
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 45457964
$ws.Range("I65").Value = 55558844
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 277794220
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -277791100
$ws.Range("N65").Value = -26240
$ws.Range("H107").Value = 404.2
$ws.Range("I107").Value = 390
$ws.Range("J107").Value = 425.5
$ws.Range("K107").Value = 1170
$ws.Range("L107").Value = 1276.5
$ws.Range("M107").Value = 750
$ws.Range("N107").Value = -5116.5
$ws.Range("H113").Value = 378.33334
$ws.Range("J113").Value = 896.5
$ws.Range("L113").Value = 2689.5
$ws.Range("N113").Value = -7029.5
$ws.Range("H122").Value = 11820751
$ws.Range("I122").Value = 14447163
$ws.Range("J122").Value = 1897.5
$ws.Range("K122").Value = 43341489
$ws.Range("L122").Value = 5692.5
$ws.Range("M122").Value = -43339039
$ws.Range("N122").Value = -10592.5
$ws.Range("H132").Value = 3737.0476
$ws.Range("I132").Value = 3126.0667
$ws.Range("J132").Value = 5264.5
$ws.Range("K132").Value = 9378.2001
$ws.Range("L132").Value = 15793.5
$ws.Range("M132").Value = -6848.2001
$ws.Range("N132").Value = -20853.5
